# Update the date on the title slide: "13/11/202" -> "22/11/202"
# The original run is split into three runs ("Date :  ", "22", "/11/202")
# by re-assigning the text of the "13" sub-range via TextRange.Characters,
# which mirrors how PowerPoint splits a run when only part of it is edited.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange

$fullText = $tr.Text
$searchText = "13/11/202"
$startIdx = $fullText.IndexOf($searchText)

if ($startIdx -ge 0) {
    # Characters() is 1-based like the rest of the PowerPoint object model.
    $target = $tr.Characters($startIdx + 1, 2)
    $target.Text = "22"
}
